$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D3"   = -7.302099999999992
    "B12"  = 5.070399999999997
    "D14"  = -7.653099999999999
    "D26"  = -8.618200000000014
    "D31"  = -8.642900000000004
    "B32"  = 6.678999999999997
    "D35"  = -7.903000000000002
    "B36"  = 9.086300000000007
    "D37"  = -7.744799999999998
    "B38"  = 5.227400000000004
    "D45"  = -7.617499999999999
    "B46"  = 6.141800000000001
    "B54"  = 4.756199999999996
    "B55"  = 5.275199999999999
    "D57"  = -8.383999999999999
    "B67"  = 5.701199999999996
    "B69"  = 5.354899999999997
    "B72"  = 5.023000000000002
    "B91"  = 5.395900000000003
    "B99"  = 4.548799999999998
    "D100" = -8.086900000000005
    "D102" = -7.766499999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
